$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed: $find"
    }
}

# p3-intro
Replace-Text ' - Our goals for project three are' ' -We choose service-oriented design because there are multiple functions, or services, that we want to provide for our user. The functions could be separately maintained and used. Our goal for project three is'

# p3-handwritings
Replace-Text 'read text and handwriting, in other words' 'read text and handwritings, in other words'

# p3-also-plan
Replace-Text 'We do plan on adding more features such as translation to this program for project four. ' 'We also plan on adding more features such as translation and pdf conversion to this program for project four, which are the other services that would be added. '

# p3-website-block
Replace-Text 'Our website supports the majority of picture file types such as pdf, png, jpg, jpeg, and gif. By having a website, it is easier for users with no computer science background to interact with our system. With the easy setup, our users should be able to navigate around our website with ease. This satisfies the black box quality of service-oriented design for our users. Since our users don’t have to be aware of out service’s inner working to use our websites to use our website. ' 'Instead of having to a long time typing out the texts in an image, the user would easily obtain the text within seconds using our program. Our website supports the majority of picture file types such as pdf, png, jpg, jpeg, and gif. By having a website, it is easier for users with no computer science background to interact with our system. The website setup would just be consisting of two boxes. One box would show the uploaded image, while the other box would show the converted text. This satisfies the black box quality of a service component for our users. Since our users don’t have to be aware of out service’s inner working to use the services on our website. With the easy setup, our users should be able to navigate around our website with ease. '

# p5-full-rewrite
Replace-Text 'So far we have two main components in our program, the logic for our OCR, and the user interface of our OCR. The logic portions of our program are in a ocr.py file, and our user interface is in a app.py file. Other than that, we have a index.html for our website and some formatting files for formatting. The program starts with the user uploading a file on our website. Then, app.py passes the image file to our ocr.py, our ocr.py reads the image and passes the read text as a string back to app.py. In our final step, our app.py displays the text on our website for the user to see. I believe our project has the pipes and filters architecture, since we out system is a form of transforming input data through a series of computational components (the image) into output data (text). It has a rather simple pipes and filters architecture, but a pipes and filter architecture, nonetheless. ' 'Our project has the pipes and filters architecture, since the system goes through a series of process of transforming input data(image) through computational components (image manipulations) into output data (text). We have two main components in our program, the logic for our OCR, and the user interface of our OCR. The logic portions of our program are in a ocr.py file, and our user interface is in a app.py file. The program starts with the user uploading a file on our website. Then, app.py passes the image file to our ocr.py. The ocr.py would then process the image in a series of image manipulations to output a cleaned image that would be used to extract the text. The ocr.py would then return the extracted text from the image and passes text as a string back to app.py. In our final step, our app.py displays the text on our website for the user to see. If we add the aspect of translation and pdf conversion in project four, it would be an extra step after processing the output text from ocr.py. For translation, the image manipulation and cleaning process would be the same. After that, we would have to determine the language that is being translated. Then the text would be translated before getting display on the website screen. For PDF conversion, the process of transformation will also be the same. There will just be an extra component in data display. The text would just be converted into a PDF that user will be able to download while also being printed on the screen. '

# p11-interpreter-to-chain
Replace-Text 'Interpreter (Behavior) – We have an image text interpreter that reads from an image then interpreted it into printed text.' 'Chain of Responsibility (Behavior) – User upload an image, app.py receive the image, pass it to ocr.py. ocr.py interoperate the text and pass it to app.py to display back to the user on our website. '

# p12-chain-to-iterator
Replace-Text 'Chain of Responsibility (Behavior) – User upload an image, app.py receive the image, pass it to ocr.py. ocr.py interoperate the text and pass it to app.py to display back to the user on our website. ' 'Iterator (Behavior) – website/display with upload and convert button etc. User will have easy access to the elements without exposing the underlying representations.'

# Delete empty paragraph between 'Use Case Diagram' paragraph and 'Identify Design Patterns' paragraph
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq [string][char]13 -and $p.Range.ListFormat.ListLevelNumber -eq 2) {
    }
}
